$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header labels for K1 / L1 ---
$ws.Range("K1").Value = "Type (for ADD)"
$ws.Range("L1").Value = "Email/Domain (for ADD)"

# --- Replace the data rows (2-11) with the updated project structure ---
$data = @(
    @("/Python_Admin_Tool_TESTING/Sub-Folder 1", "Sub-Folder 1", "1jPqgww8lNGleK7h15iHuNdyUbZDf0idz", "Commenter", "domain", "bioaccessla.com", "jmoreno@bioaccessla.com", "https://drive.google.com/drive/folders/1jPqgww8lNGleK7h15iHuNdyUbZDf0idz"),
    @("/Python_Admin_Tool_TESTING/Sub-Folder 1", "Sub-Folder 1", "1jPqgww8lNGleK7h15iHuNdyUbZDf0idz", "Viewer", "user", "ernie.moreno62@gmail.com", "jmoreno@bioaccessla.com", "https://drive.google.com/drive/folders/1jPqgww8lNGleK7h15iHuNdyUbZDf0idz"),
    @("/Python_Admin_Tool_TESTING/Sub-Folder 1", "Sub-Folder 1", "1jPqgww8lNGleK7h15iHuNdyUbZDf0idz", "Owner", "user", "jmoreno@bioaccessla.com", "jmoreno@bioaccessla.com", "https://drive.google.com/drive/folders/1jPqgww8lNGleK7h15iHuNdyUbZDf0idz"),
    @("/Python_Admin_Tool_TESTING/Sub-Folder 1/Test Sheet 2", "Test Sheet 2", "1Wan1C_Cxndc2M6yXKa8vxJLXCgYoBXFzDVwRYAfGbVY", "Editor", "group", "jesus_test_group@bioaccessla.com", "jmoreno@bioaccessla.com", "https://docs.google.com/spreadsheets/d/1Wan1C_Cxndc2M6yXKa8vxJLXCgYoBXFzDVwRYAfGbVY/edit?usp=drivesdk"),
    @("/Python_Admin_Tool_TESTING/Sub-Folder 1/Test Sheet 2", "Test Sheet 2", "1Wan1C_Cxndc2M6yXKa8vxJLXCgYoBXFzDVwRYAfGbVY", "Commenter", "domain", "bioaccessla.com", "jmoreno@bioaccessla.com", "https://docs.google.com/spreadsheets/d/1Wan1C_Cxndc2M6yXKa8vxJLXCgYoBXFzDVwRYAfGbVY/edit?usp=drivesdk"),
    @("/Python_Admin_Tool_TESTING/Sub-Folder 1/Test Sheet 2", "Test Sheet 2", "1Wan1C_Cxndc2M6yXKa8vxJLXCgYoBXFzDVwRYAfGbVY", "Viewer", "user", "ernie.moreno62@gmail.com", "jmoreno@bioaccessla.com", "https://docs.google.com/spreadsheets/d/1Wan1C_Cxndc2M6yXKa8vxJLXCgYoBXFzDVwRYAfGbVY/edit?usp=drivesdk"),
    @("/Python_Admin_Tool_TESTING/Sub-Folder 1/Test Sheet 2", "Test Sheet 2", "1Wan1C_Cxndc2M6yXKa8vxJLXCgYoBXFzDVwRYAfGbVY", "Owner", "user", "jmoreno@bioaccessla.com", "jmoreno@bioaccessla.com", "https://docs.google.com/spreadsheets/d/1Wan1C_Cxndc2M6yXKa8vxJLXCgYoBXFzDVwRYAfGbVY/edit?usp=drivesdk"),
    @("/Python_Admin_Tool_TESTING/Test Doc 1", "Test Doc 1", "1O90b5jSuK3lIz-RYZIEtAAlA3-IQ_vmxgulyB_6vY2U", "Viewer", "domain", "bioaccessla.com", "jmoreno@bioaccessla.com", "https://docs.google.com/document/d/1O90b5jSuK3lIz-RYZIEtAAlA3-IQ_vmxgulyB_6vY2U/edit?usp=drivesdk"),
    @("/Python_Admin_Tool_TESTING/Test Doc 1", "Test Doc 1", "1O90b5jSuK3lIz-RYZIEtAAlA3-IQ_vmxgulyB_6vY2U", "Editor", "user", "ernie.moreno62@gmail.com", "jmoreno@bioaccessla.com", "https://docs.google.com/document/d/1O90b5jSuK3lIz-RYZIEtAAlA3-IQ_vmxgulyB_6vY2U/edit?usp=drivesdk"),
    @("/Python_Admin_Tool_TESTING/Test Doc 1", "Test Doc 1", "1O90b5jSuK3lIz-RYZIEtAAlA3-IQ_vmxgulyB_6vY2U", "Owner", "user", "jmoreno@bioaccessla.com", "jmoreno@bioaccessla.com", "https://docs.google.com/document/d/1O90b5jSuK3lIz-RYZIEtAAlA3-IQ_vmxgulyB_6vY2U/edit?usp=drivesdk")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# --- Conditional formatting over A2:K1048576 ---
# Rules: ADD -> light green, REMOVE -> light red, MODIFY -> light yellow
$cfRange = $ws.Range("A2:K1048576")
$cfRange.FormatConditions.Delete()

$fc1 = $cfRange.FormatConditions.Add(2, 3, '=$I2="ADD"')
$fc1.Interior.Color = 12315096

$fc2 = $cfRange.FormatConditions.Add(2, 3, '=$I2="REMOVE"')
$fc2.Interior.Color = 13551615

$fc3 = $cfRange.FormatConditions.Add(2, 3, '=$I2="MODIFY"')
$fc3.Interior.Color = 10284031
